$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 32500
$ws.Range("I21").Value = 36666.668
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 36666.668
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -36198.668
$ws.Range("N21").Value = -20936

$ws.Range("H23").Value = 32500
$ws.Range("I23").Value = 36666.668
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 36666.668
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -36432.668
$ws.Range("N23").Value = -20468

$ws.Range("H86").Value = 2968.8667
$ws.Range("I86").Value = 2733.3333
$ws.Range("J86").Value = 3322.1667
$ws.Range("K86").Value = 2733.3333
$ws.Range("L86").Value = 3322.1667
$ws.Range("M86").Value = -1610.3333
$ws.Range("N86").Value = -5568.1667

$ws.Range("H89").Value = 2968.8667
$ws.Range("I89").Value = 2733.3333
$ws.Range("J89").Value = 3322.1667
$ws.Range("K89").Value = 13666.6665
$ws.Range("L89").Value = 16610.8335
$ws.Range("M89").Value = -8050.666499999999
$ws.Range("N89").Value = -27842.8335

$ws.Range("H129").Value = 1377.7693
$ws.Range("I129").Value = 663.5
$ws.Range("J129").Value = 1507.6364
$ws.Range("K129").Value = 1990.5
$ws.Range("L129").Value = 4522.9092
$ws.Range("M129").Value = 3009.5
$ws.Range("N129").Value = -14522.9092

$ws.Range("H131").Value = 8110.952
$ws.Range("I131").Value = 1886
$ws.Range("J131").Value = 10056.25
$ws.Range("K131").Value = 5658
$ws.Range("L131").Value = 30168.75
$ws.Range("M131").Value = -618
$ws.Range("N131").Value = -40248.75

$ws.Range("H135").Value = 552.94446
$ws.Range("I135").Value = 479.2456
$ws.Range("J135").Value = 833
$ws.Range("K135").Value = 4313.2104
$ws.Range("L135").Value = 7497
$ws.Range("M135").Value = -1778.2104
$ws.Range("N135").Value = -12567

$ws.Range("H137").Value = 1217.3055
$ws.Range("I137").Value = 1170.3939
$ws.Range("K137").Value = 3511.1817
$ws.Range("M137").Value = -961.1817000000001

$ws.Range("H138").Value = 2051.6619
$ws.Range("J138").Value = 2611.0881
$ws.Range("L138").Value = 7833.2643
$ws.Range("N138").Value = -18113.2643

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 50000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H25").Value = 15003.6
$ws.Range("I25").Value = 13250
$ws.Range("J25").Value = 22018
$ws.Range("K25").Value = 13250
$ws.Range("L25").Value = 22018
$ws.Range("M25").Value = -12848
$ws.Range("N25").Value = -22822

$ws.Range("H88").Value = 4123.3335
$ws.Range("I88").Value = 2685
$ws.Range("J88").Value = 7000
$ws.Range("K88").Value = 2685
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = -2279
$ws.Range("N88").Value = -7812

$ws.Range("H91").Value = 4123.3335
$ws.Range("I91").Value = 2685
$ws.Range("J91").Value = 7000
$ws.Range("K91").Value = 2685
$ws.Range("L91").Value = 7000
$ws.Range("M91").Value = -1281
$ws.Range("N91").Value = -9808

$ws.Range("H122").Value = 4335.7617
$ws.Range("I122").Value = 4134.1665
$ws.Range("J122").Value = 4604.5557
$ws.Range("K122").Value = 12402.4995
$ws.Range("L122").Value = 13813.6671
$ws.Range("M122").Value = -9952.499500000002
$ws.Range("N122").Value = -18713.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33319.188
$ws.Range("I20").Value = 60358.47
$ws.Range("J20").Value = 2674.6667
$ws.Range("K20").Value = 60358.47
$ws.Range("L20").Value = 2674.6667
$ws.Range("M20").Value = -60111.47
$ws.Range("N20").Value = -3168.6667

$ws.Range("H107").Value = 14501.581
$ws.Range("I107").Value = 27848.4
$ws.Range("K107").Value = 27848.4
$ws.Range("M107").Value = -25928.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 20004
$ws.Range("J6").Value = 20004
$ws.Range("L6").Value = 20004
$ws.Range("N6").Value = -20230

$ws.Range("H16").Value = 20004
$ws.Range("J16").Value = 20004
$ws.Range("L16").Value = 20004
$ws.Range("N16").Value = -20504

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H38").Value = 41800
$ws.Range("J38").Value = 39750
$ws.Range("L38").Value = 39750
$ws.Range("N38").Value = -40676

$ws.Range("H40").Value = 27900
$ws.Range("J40").Value = 5800
$ws.Range("L40").Value = 5800
$ws.Range("N40").Value = -6102

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H46").Value = 15950
$ws.Range("J46").Value = 19933.334
$ws.Range("L46").Value = 19933.334
$ws.Range("N46").Value = -20245.334

$ws.Range("H52").Value = 500000
$ws.Range("J52").Value = 500000
$ws.Range("L52").Value = 500000
$ws.Range("N52").Value = -500518

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()

$ws.Range("H58").Value = 50000
$ws.Range("J58").Value = 50000
$ws.Range("L58").Value = 50000
$ws.Range("N58").Value = -50554

$ws.Range("H102").Value = 2152.3547
$ws.Range("I102").Value = 1939.0952
$ws.Range("K102").Value = 1939.0952
$ws.Range("M102").Value = -317.0952

$ws.Range("H126").Value = 2824.5
$ws.Range("I126").Value = 1899
$ws.Range("K126").Value = 5697
$ws.Range("M126").Value = -3227

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3998.5
$ws.Range("J3").Value = 3998.5
$ws.Range("L3").Value = 3998.5
$ws.Range("N3").Value = -4222.5

$ws.Range("H15").Value = 3998.5
$ws.Range("J15").Value = 3998.5
$ws.Range("L15").Value = 3998.5
$ws.Range("N15").Value = -4338.5

$ws.Range("H40").Value = 4425.8
$ws.Range("I40").Value = 4657.25
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 4657.25
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -4521.25
$ws.Range("N40").Value = -3772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 30000
$ws.Range("J31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("N31").Value = -30696

$ws.Range("H126").Value = 4767.952
$ws.Range("J126").Value = 921.8333
$ws.Range("L126").Value = 2765.4999
$ws.Range("N126").Value = -7705.4999
